# Type Effectiveness Chart update
# - fixes several effectiveness values in the C3:I9 grid
# - adds a per-row total column (J) with SUM formulas
# - adds a per-column total row (10) with SUM formulas
# - nudges the alignment on a few cells that flip to their "weak" (0.5) value
# - updates the saved selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Correct the effectiveness values that differ from the original chart.
#    (row => column letter => new value)
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 9).Value  = 0.5   # I3
$ws.Cells.Item(4, 4).Value  = 0.5   # D4
$ws.Cells.Item(4, 9).Value  = 1     # I4
$ws.Cells.Item(5, 5).Value  = 0.5   # E5
$ws.Cells.Item(5, 9).Value  = 2     # I5
$ws.Cells.Item(6, 6).Value  = 0.5   # F6
$ws.Cells.Item(6, 9).Value  = 2     # I6
$ws.Cells.Item(7, 9).Value  = 0.5   # I7
$ws.Cells.Item(8, 8).Value  = 0.5   # H8
$ws.Cells.Item(8, 9).Value  = 1     # I8
$ws.Cells.Item(9, 3).Value  = 2     # C9
$ws.Cells.Item(9, 4).Value  = 1     # D9
$ws.Cells.Item(9, 5).Value  = 0.5   # E9
$ws.Cells.Item(9, 6).Value  = 0.5   # F9
$ws.Cells.Item(9, 7).Value  = 2     # G9
$ws.Cells.Item(9, 8).Value  = 1     # H9
$ws.Cells.Item(9, 9).Value  = 0.5   # I9

# A handful of cells that now hold a "resisted" (0.5) value also pick up an
# explicit center alignment (re-applying the same alignment the rest of the
# sheet already uses).
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("H8").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 2. Row totals in column J (J3:J9), one SUM per attacking-type row.
#    J3 is entered on its own; J4:J9 is entered as a single fill so the
#    engine keeps it as one relative formula across the block.
# ---------------------------------------------------------------------------
$ws.Range("J3").Formula = "=SUM(C3:I3)"
$ws.Range("J4:J9").Formula = "=SUM(C4:I4)"

# ---------------------------------------------------------------------------
# 3. Column totals in row 10 (C10:I10), one SUM per defending-type column.
#    C10 is entered on its own; D10:I10 is entered as a single fill.
# ---------------------------------------------------------------------------
$ws.Range("C10").Formula = "=SUM(C3:C9)"
$ws.Range("D10:I10").Formula = "=SUM(D3:D9)"

# ---------------------------------------------------------------------------
# 4. Restore the last-used selection that was saved with the workbook.
# ---------------------------------------------------------------------------
$ws.Range("N9").Select()

Write-Output "Type Effectiveness Chart updated"
